# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 7490
$ws1.Range("F4").Value = 3560
$ws1.Range("F6").Value = 3877
$ws1.Range("F18").Value = 4214
$ws1.Range("F21").Value = 1035
$ws1.Range("F23").Value = 1916
$ws1.Range("F37").Value = 4440
$ws1.Range("F39").Value = 330
$ws1.Range("F42").Value = 846

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7490
$ws4.Range("F5").Value = 3560
$ws4.Range("F6").Value = 3877
$ws4.Range("F19").Value = 4214
$ws4.Range("F26").Value = 1916
$ws4.Range("F36").Value = 4440
$ws4.Range("F39").Value = 330
$ws4.Range("F42").Value = 846

$wb.Save()
